# Update Leve profit calculation columns (H-N) across multiple sheets
# to reflect refreshed market-board pricing data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1157.2222
$ws.Range("I18").Value = 932.1429000000001
$ws.Range("J18").Value = 1945
$ws.Range("K18").Value = 932.1429000000001
$ws.Range("L18").Value = 1945
$ws.Range("M18").Value = -648.1429000000001
$ws.Range("N18").Value = -2513

$ws.Range("H33").Value = 436.52173
$ws.Range("I33").Value = 420.5
$ws.Range("K33").Value = 420.5
$ws.Range("M33").Value = -191.5

$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()

$ws.Range("H76").Value = 5664.1665
$ws.Range("I76").Value = 6000
$ws.Range("K76").Value = 6000
$ws.Range("M76").Value = -5685

$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()

$ws.Range("H79").Value = 5664.1665
$ws.Range("I79").Value = 6000
$ws.Range("K79").Value = 6000
$ws.Range("M79").Value = -4908

$ws.Range("H106").Value = 8595.25
$ws.Range("I106").Value = 8968.267
$ws.Range("K106").Value = 8968.267
$ws.Range("M106").Value = -8337.267

$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()

$ws.Range("H132").Value = 8338883
$ws.Range("I132").Value = 12352285
$ws.Range("J132").Value = 3355.1538
$ws.Range("K132").Value = 37056855
$ws.Range("L132").Value = 10065.4614
$ws.Range("M132").Value = -37054325
$ws.Range("N132").Value = -15125.4614

$ws.Range("H137").Value = 1809.3529
$ws.Range("J137").Value = 2105.8462
$ws.Range("L137").Value = 6317.5386
$ws.Range("N137").Value = -11417.5386

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 62501320
$ws.Range("I61").Value = 76924160
$ws.Range("K61").Value = 76924160
$ws.Range("M61").Value = -76923948

$ws.Range("H74").Value = 2957.0908
$ws.Range("I74").Value = 1950
$ws.Range("J74").Value = 3180.889
$ws.Range("K74").Value = 1950
$ws.Range("L74").Value = 3180.889
$ws.Range("M74").Value = -1076
$ws.Range("N74").Value = -4928.889

$ws.Range("H77").Value = 2957.0908
$ws.Range("I77").Value = 1950
$ws.Range("J77").Value = 3180.889
$ws.Range("K77").Value = 9750
$ws.Range("L77").Value = 15904.445
$ws.Range("M77").Value = -5382
$ws.Range("N77").Value = -24640.445

$ws.Range("H97").Value = 388.86957
$ws.Range("I97").Value = 352.2
$ws.Range("K97").Value = 352.2
$ws.Range("M97").Value = 143.8

$ws.Range("I122").Value = 1559.75
$ws.Range("J122").Value = 2118.1428
$ws.Range("K122").Value = 4679.25
$ws.Range("L122").Value = 6354.428400000001
$ws.Range("M122").Value = -2229.25
$ws.Range("N122").Value = -11254.4284

$ws.Range("H132").Value = 2368.319
$ws.Range("I132").Value = 2259.0454
$ws.Range("J132").Value = 2464.48
$ws.Range("K132").Value = 6777.1362
$ws.Range("L132").Value = 7393.440000000001
$ws.Range("M132").Value = -4247.1362
$ws.Range("N132").Value = -12453.44

$ws.Range("H136").Value = 62501320
$ws.Range("I136").Value = 76924160
$ws.Range("K136").Value = 230772480
$ws.Range("M136").Value = -230769930

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1899.75
$ws.Range("I134").Value = 1547.3334
$ws.Range("K134").Value = 4642.0002
$ws.Range("M134").Value = -2107.0002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 44102.125
$ws.Range("I22").Value = 212.2
$ws.Range("K22").Value = 212.2
$ws.Range("M22").Value = 137.8

$ws.Range("H99").Value = 1616.8889
$ws.Range("I99").Value = 1594
$ws.Range("K99").Value = 1594
$ws.Range("M99").Value = -96

$ws.Range("H105").Value = 1044.2
$ws.Range("I105").Value = 1070
$ws.Range("J105").Value = 1005.5
$ws.Range("K105").Value = 1070
$ws.Range("L105").Value = 1005.5
$ws.Range("M105").Value = 677
$ws.Range("N105").Value = -4499.5

$ws.Range("H122").Value = 1093.75
$ws.Range("I122").Value = 1093.75
$ws.Range("K122").Value = 3281.25
$ws.Range("M122").Value = -831.25

$ws.Range("H126").Value = 1616.8889
$ws.Range("I126").Value = 1594
$ws.Range("K126").Value = 4782
$ws.Range("M126").Value = -2312

$ws.Range("H134").Value = 14287348
$ws.Range("I134").Value = 1421.0435
$ws.Range("K134").Value = 4263.1305
$ws.Range("M134").Value = -1728.1305

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 729.06976
$ws.Range("I5").Value = 730.375
$ws.Range("K5").Value = 2191.125
$ws.Range("M5").Value = -2079.125

$ws.Range("H68").Value = 905.6429000000001
$ws.Range("I68").Value = 1300
$ws.Range("J68").Value = 798.0909
$ws.Range("K68").Value = 3900
$ws.Range("L68").Value = 2394.2727
$ws.Range("M68").Value = -3089
$ws.Range("N68").Value = -4016.2727

$ws.Range("H71").Value = 905.6429000000001
$ws.Range("I71").Value = 1300
$ws.Range("J71").Value = 798.0909
$ws.Range("K71").Value = 11700
$ws.Range("L71").Value = 7182.8181
$ws.Range("M71").Value = -7644
$ws.Range("N71").Value = -15294.8181

$ws.Range("H131").Value = 15874435
$ws.Range("J131").Value = 1528.9286
$ws.Range("L131").Value = 4586.7858
$ws.Range("N131").Value = -14666.7858

$ws.Range("H135").Value = 729.06976
$ws.Range("I135").Value = 730.375
$ws.Range("K135").Value = 6573.375
$ws.Range("M135").Value = -4038.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 50003172
$ws.Range("I70").Value = 35717804
$ws.Range("J70").Value = 100001950
$ws.Range("K70").Value = 35717804
$ws.Range("L70").Value = 100001950
$ws.Range("M70").Value = -35717534
$ws.Range("N70").Value = -100002490

$ws.Range("H73").Value = 50003172
$ws.Range("I73").Value = 35717804
$ws.Range("J73").Value = 100001950
$ws.Range("K73").Value = 35717804
$ws.Range("L73").Value = 100001950
$ws.Range("M73").Value = -35716868
$ws.Range("N73").Value = -100003822

$ws.Range("H102").Value = 3178.5
$ws.Range("I102").Value = 2600
$ws.Range("J102").Value = 3757
$ws.Range("K102").Value = 2600
$ws.Range("L102").Value = 3757
$ws.Range("M102").Value = -978
$ws.Range("N102").Value = -7001

$ws.Range("H109").Value = 6856.5713
$ws.Range("J109").Value = 6856.5713
$ws.Range("L109").Value = 6856.5713
$ws.Range("N109").Value = -8936.5713

$ws.Range("H113").Value = 1492.9231
$ws.Range("I113").Value = 1450.8889
$ws.Range("J113").Value = 1587.5
$ws.Range("K113").Value = 1450.8889
$ws.Range("L113").Value = 1587.5
$ws.Range("M113").Value = 719.1111000000001
$ws.Range("N113").Value = -5927.5

$ws.Range("H132").Value = 2670.6667
$ws.Range("I132").Value = 2323.8462
$ws.Range("K132").Value = 6971.5386
$ws.Range("M132").Value = -4441.5386

$ws.Range("H135").Value = 34340.715
$ws.Range("J135").Value = 33760.742
$ws.Range("L135").Value = 33760.742
$ws.Range("N135").Value = -43900.742

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6256.2
$ws.Range("I40").Value = 2575.5
$ws.Range("K40").Value = 2575.5
$ws.Range("M40").Value = -2439.5

$ws.Range("H122").Value = 62501500
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

$ws.Range("H136").Value = 1905.4286
$ws.Range("I136").Value = 1467.8
$ws.Range("K136").Value = 4403.4
$ws.Range("M136").Value = -1853.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 408.0625
$ws.Range("I113").Value = 335.22223
$ws.Range("K113").Value = 1005.66669
$ws.Range("M113").Value = 1164.33331

$ws.Range("H122").Value = 16669202
$ws.Range("I122").Value = 17859716
$ws.Range("J122").Value = 2005
$ws.Range("K122").Value = 53579148
$ws.Range("L122").Value = 6015
$ws.Range("M122").Value = -53576698
$ws.Range("N122").Value = -10915

$ws.Range("H126").Value = 35715044
$ws.Range("I126").Value = 50000430
$ws.Range("J126").Value = 1577.75
$ws.Range("K126").Value = 150001290
$ws.Range("L126").Value = 4733.25
$ws.Range("M126").Value = -149998820
$ws.Range("N126").Value = -9673.25

$ws.Range("H136").Value = 1689.2632
$ws.Range("I136").Value = 1509.7333
$ws.Range("J136").Value = 2362.5
$ws.Range("K136").Value = 4529.199900000001
$ws.Range("L136").Value = 7087.5
$ws.Range("M136").Value = -1979.199900000001
$ws.Range("N136").Value = -12187.5
